$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the Date column (A2:A27) contents while keeping the existing cell
# formatting/style intact (matches the workbook change where the <v> date
# serials were removed but the style attribute s="1" was kept).
$ws.Range("A2:A27").ClearContents()

# Update Activity (column D) and Minutes (column F) values for the rows
# that changed.
$ws.Range("D4").Value = "Online module"
$ws.Range("F4").Value = 240

$ws.Range("D5").Value = "Computer exercises"

$ws.Range("F7").Value = 240

$ws.Range("D8").Value = "Online module"
$ws.Range("F8").Value = 240

$ws.Range("D9").Value = "Computer exercises"

$ws.Range("D10").Value = "Online module"
$ws.Range("F10").Value = 240

$ws.Range("F11").Value = 240

$ws.Range("D12").Value = "Online module"
$ws.Range("F12").Value = 240

$ws.Range("F13").Value = 240

$ws.Range("D15").Value = "Online module"
$ws.Range("F15").Value = 240

$ws.Range("F16").Value = 240

$ws.Range("D18").Value = "Computer exercises"
$ws.Range("F18").Value = 240

$ws.Range("F19").Value = 240

$ws.Range("D20").Value = "Online module"
$ws.Range("F20").Value = 240

$ws.Range("D21").Value = "Computer exercises"

$ws.Range("F24").Value = 120

$ws.Range("D26").Value = "Online module"
$ws.Range("F26").Value = 240

# Update the active selection to match the saved view state.
$ws.Range("D22").Select()
